$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kotte2014")

# Add the new row (row 9) mirroring row 8's pattern, with new enzyme "ENZtr"
$ws.Range("A9").Value = "ENZtr"
$ws.Range("C9").Value = "enz[c] <==>"
$ws.Range("E9").Value = 0
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1

# Update the selection to match the committed view state
$ws.Range("C8:C9").Select()
